# Update "想去人数" (interested-count) values in column F on the
# "展览" and "全部类型" worksheets, per the source diff.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 352
$ws1.Range("F7").Value = 898
$ws1.Range("F8").Value = 67
$ws1.Range("F9").Value = 547
$ws1.Range("F15").Value = 45
$ws1.Range("F17").Value = 6729
$ws1.Range("F19").Value = 77
$ws1.Range("F20").Value = 24
$ws1.Range("F21").Value = 7652
$ws1.Range("F24").Value = 3421
$ws1.Range("F26").Value = 2154
$ws1.Range("F27").Value = 920
$ws1.Range("F29").Value = 186
$ws1.Range("F35").Value = 1783
$ws1.Range("F39").Value = 8
$ws1.Range("F41").Value = 1251
$ws1.Range("F42").Value = 1884

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 352
$ws4.Range("F9").Value = 898
$ws4.Range("F10").Value = 67
$ws4.Range("F11").Value = 547
$ws4.Range("F18").Value = 45
$ws4.Range("F20").Value = 6729
$ws4.Range("F22").Value = 77
$ws4.Range("F23").Value = 24
$ws4.Range("F24").Value = 7652
$ws4.Range("F27").Value = 3421
$ws4.Range("F29").Value = 2154
$ws4.Range("F30").Value = 920
$ws4.Range("F32").Value = 186
$ws4.Range("F38").Value = 1783
$ws4.Range("F42").Value = 8
$ws4.Range("F44").Value = 1251
$ws4.Range("F45").Value = 1884
